$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59; this shifts existing rows 59..81 down to 60..82
$ws.Rows(59).Insert()

# Populate the new row 59 with the new data entry
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44845
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = 100112012
$ws.Cells.Item(59, 7).Value = "Espinaca"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 250
$ws.Cells.Item(59, 11).Value = 1800
$ws.Cells.Item(59, 12).Value = 2000
$ws.Cells.Item(59, 13).Value = 1900
$ws.Cells.Item(59, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(59, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value = 633
$ws.Cells.Item(59, 17).Value = 3
$ws.Cells.Item(59, 18).Value = "Hortaliza"
